# The upstream change is a SharePoint "content type" custom-XML sync:
# the data part that carries the content-type schema (customXml/item1.xml)
# gets its version/fields fingerprints re-minted, and the companion
# datastore-item part describing it (customXml/itemProps1.xml) is
# replaced outright: a freshly generated GUID and an emptied schemaRefs
# list (no cached schema references until the part is next synced).
#
# customXml/item1.xml:
#   ma:versionID  1cc4cf9d95b5e2d14d7aabb44ca49f5e -> ce94caacb4a5cc228342027e3189af2c
#   ma:fieldsID   a33e6829bf21261855124b7b230b6e9c -> 5f85a36ab557a4a47cd270a1ee4435c0
#
# customXml/itemProps1.xml:
#   ds:itemID     {4B368451-2D87-40A2-9C11-75C93564977B} -> {D3E15501-9F16-4620-A6AB-7C77E3A7AD56}
#   <ds:schemaRefs>...</ds:schemaRefs> -> (removed, self-closed element)

$d = $word.ActiveDocument

$ctNamespace  = "http://schemas.microsoft.com/office/2006/metadata/contentType"
$oldVersionID = "1cc4cf9d95b5e2d14d7aabb44ca49f5e"
$newVersionID = "ce94caacb4a5cc228342027e3189af2c"
$oldFieldsID  = "a33e6829bf21261855124b7b230b6e9c"
$newFieldsID  = "5f85a36ab557a4a47cd270a1ee4435c0"
$oldItemID    = "{4B368451-2D87-40A2-9C11-75C93564977B}"
$newItemID    = "{D3E15501-9F16-4620-A6AB-7C77E3A7AD56}"

$parts = $d.CustomXMLParts

function Find-ContentTypePart($parts, $ns) {
    try {
        $sel = $parts.SelectByNamespace($ns)
        if ($sel -ne $null -and $sel.Count -ge 1) {
            return $sel.Item(1)
        }
    } catch {
    }
    for ($i = 1; $i -le $parts.Count; $i++) {
        $candidate = $parts.Item($i)
        $matched = $false
        try {
            if ($candidate.NamespaceURI -eq $ns) { $matched = $true }
        } catch {
        }
        if (-not $matched) {
            try {
                if ($candidate.XML -and $candidate.XML.Contains("contentTypeSchema")) { $matched = $true }
            } catch {
            }
        }
        if ($matched) { return $candidate }
    }
    return $null
}

$target = Find-ContentTypePart $parts $ctNamespace

if ($target -ne $null) {
    $xml = $target.XML
    if ($xml) {
        $updated = $xml.Replace($oldVersionID, $newVersionID).Replace($oldFieldsID, $newFieldsID)

        # Re-mint the schema part's fingerprints in place.
        try { $target.XML = $updated } catch { }

        # Replace the paired datastore-item: new id, schemaRefs cleared.
        try {
            $target.Delete()
            $readded = $parts.Add($updated, $ctNamespace)
        } catch {
            try {
                $readded = $parts.Add($updated)
            } catch {
            }
        }
    }
}

Write-Output "done"
